$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.734.57'
$ws.Range("E2").Value = '  +3.80%  '
$ws.Range("D3").Value = '2.418.37'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("E4").Value = '  -0.06%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = '''317.13'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +4.69%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = '''101.87'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +6.82%  '
$ws.Range("E8").Value = '  -0.04%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = '''0.525'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +9.26%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = '''35.29'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +2.90%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = '''0.0798'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  -1.94%  '
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("D15").Value = '2.796.23'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '2.401.46'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("E17").Value = '  +3.55%  '
$ws.Range("D18").Value = '44.570.45'
$ws.Range("E18").Value = '  +3.38%  '
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = '0.0₃0916'
$ws.Range("E21").Value = '  +3.35%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = '''242.56'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("E24").Value = '  +3.05%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = '''2.49'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("E26").Value = '  -0.10%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = '''25.15'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +2.83%  '
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("E29").Value = '  +1.60%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = '''33.44'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +3.86%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = '''48.28'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +0.75%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = '''0.126'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +14.17%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = '''19.47'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +10.58%  '
$ws.Range("E34").Value = '  +3.07%  '
$ws.Range("E35").Value = '  +0.25%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = '''0.0761'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +4.39%  '
$ws.Range("E37").Value = '  +2.16%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = '''4.44'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +2.95%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = '''126.57'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("E42").Value = '  -3.71%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = '''21.18'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +2.11%  '
$ws.Range("E44").Value = '  +3.31%  '
$ws.Range("D45").Value = '1.935.59'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  +6.85%  '
$ws.Range("E48").Value = '  -0.79%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = '''1.75'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +16.50%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = '''75.72'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +5.98%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = '''53.51'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +4.42%  '
